$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5000
$ws.Range("I43").Value = 8100
$ws.Range("J43").Value = 1900
$ws.Range("K43").Value = 8100
$ws.Range("L43").Value = 1900
$ws.Range("M43").Value = -8031
$ws.Range("N43").Value = -2038
$ws.Range("H123").Value = 26164.285
$ws.Range("J123").Value = 26164.285
$ws.Range("L123").Value = 26164.285
$ws.Range("N123").Value = -35964.285
$ws.Range("H137").Value = 910.6286
$ws.Range("I137").Value = 724
$ws.Range("J137").Value = 1657.1428
$ws.Range("K137").Value = 2172
$ws.Range("L137").Value = 4971.428400000001
$ws.Range("M137").Value = 378
$ws.Range("N137").Value = -10071.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1681.8182
$ws.Range("I63").Value = 1681.8182
$ws.Range("K63").Value = 1681.8182
$ws.Range("M63").Value = -995.8181999999999
$ws.Range("H66").Value = 1681.8182
$ws.Range("I66").Value = 1681.8182
$ws.Range("K66").Value = 8409.091
$ws.Range("M66").Value = -4977.091
$ws.Range("H74").Value = 905.0833
$ws.Range("I74").Value = 1076.8334
$ws.Range("K74").Value = 1076.8334
$ws.Range("M74").Value = -202.8334
$ws.Range("H77").Value = 905.0833
$ws.Range("I77").Value = 1076.8334
$ws.Range("K77").Value = 5384.166999999999
$ws.Range("M77").Value = -1016.166999999999
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 179.625
$ws.Range("I22").Value = 191.84616
$ws.Range("J22").Value = 126.666664
$ws.Range("K22").Value = 191.84616
$ws.Range("L22").Value = 126.666664
$ws.Range("M22").Value = -18.84616
$ws.Range("N22").Value = -472.666664
$ws.Range("H64").Value = 967.5833
$ws.Range("I64").Value = 975
$ws.Range("J64").Value = 963.875
$ws.Range("K64").Value = 975
$ws.Range("L64").Value = 963.875
$ws.Range("M64").Value = -750
$ws.Range("N64").Value = -1413.875
$ws.Range("H67").Value = 967.5833
$ws.Range("I67").Value = 975
$ws.Range("J67").Value = 963.875
$ws.Range("K67").Value = 975
$ws.Range("L67").Value = 963.875
$ws.Range("M67").Value = -195
$ws.Range("N67").Value = -2523.875
$ws.Range("H81").Value = 20150
$ws.Range("J81").Value = 20150
$ws.Range("L81").Value = 20150
$ws.Range("N81").Value = -22272
$ws.Range("H84").Value = 20150
$ws.Range("J84").Value = 20150
$ws.Range("L84").Value = 60450
$ws.Range("N84").Value = -71058
$ws.Range("H86").Value = 1898.8334
$ws.Range("I86").Value = 1948.4
$ws.Range("J86").Value = 1799.7
$ws.Range("K86").Value = 1948.4
$ws.Range("L86").Value = 1799.7
$ws.Range("M86").Value = -825.4000000000001
$ws.Range("N86").Value = -4045.7
$ws.Range("H89").Value = 1898.8334
$ws.Range("I89").Value = 1948.4
$ws.Range("J89").Value = 1799.7
$ws.Range("K89").Value = 9742
$ws.Range("L89").Value = 8998.5
$ws.Range("M89").Value = -4126
$ws.Range("N89").Value = -20230.5
$ws.Range("H94").Value = 992.4783
$ws.Range("I94").Value = 948.2308
$ws.Range("J94").Value = 1050
$ws.Range("K94").Value = 948.2308
$ws.Range("L94").Value = 1050
$ws.Range("M94").Value = -497.2308
$ws.Range("N94").Value = -1952
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 23129.9
$ws.Range("I86").Value = 42499.8
$ws.Range("J86").Value = 3760
$ws.Range("K86").Value = 42499.8
$ws.Range("L86").Value = 3760
$ws.Range("M86").Value = -41376.8
$ws.Range("N86").Value = -6006
$ws.Range("H89").Value = 23129.9
$ws.Range("I89").Value = 42499.8
$ws.Range("J89").Value = 3760
$ws.Range("K89").Value = 212499
$ws.Range("L89").Value = 18800
$ws.Range("M89").Value = -206883
$ws.Range("N89").Value = -30032
$ws.Range("H99").Value = 2179.689
$ws.Range("I99").Value = 1710.4615
$ws.Range("K99").Value = 1710.4615
$ws.Range("M99").Value = -212.4614999999999
$ws.Range("H107").Value = 285.03845
$ws.Range("I107").Value = 225.21053
$ws.Range("J107").Value = 319.48486
$ws.Range("K107").Value = 225.21053
$ws.Range("L107").Value = 319.48486
$ws.Range("M107").Value = 1694.78947
$ws.Range("N107").Value = -4159.48486
$ws.Range("H126").Value = 2179.689
$ws.Range("I126").Value = 1710.4615
$ws.Range("K126").Value = 5131.3845
$ws.Range("M126").Value = -2661.3845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 579.9286
$ws.Range("I113").Value = 565
$ws.Range("J113").Value = 582.4167
$ws.Range("K113").Value = 1695
$ws.Range("L113").Value = 1747.2501
$ws.Range("M113").Value = 475
$ws.Range("N113").Value = -6087.2501
$ws.Range("H122").Value = 1668667.4
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1668667.4
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 15018006.6
$ws.Range("N122").Value = -15022906.6
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5675.8
$ws.Range("I70").Value = 4185
$ws.Range("K70").Value = 4185
$ws.Range("M70").Value = -3915
$ws.Range("H73").Value = 5675.8
$ws.Range("I73").Value = 4185
$ws.Range("K73").Value = 4185
$ws.Range("M73").Value = -3249
$ws.Range("H113").Value = 6090
$ws.Range("I113").Value = 6701.25
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 6701.25
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = -4531.25
$ws.Range("N113").Value = -5540
$ws.Range("H126").Value = 1808.3
$ws.Range("I126").Value = 1979.2727
$ws.Range("K126").Value = 5937.8181
$ws.Range("M126").Value = -3467.8181
$ws.Range("H132").Value = 2806.276
$ws.Range("I132").Value = 2437.4285
$ws.Range("K132").Value = 7312.2855
$ws.Range("M132").Value = -4782.2855
$ws.Range("H134").Value = 19999.334
$ws.Range("J134").Value = 19999.334
$ws.Range("L134").Value = 59998.00199999999
$ws.Range("N134").Value = -65068.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1047.1428
$ws.Range("I46").Value = 1271.1
$ws.Range("J46").Value = 487.25
$ws.Range("K46").Value = 1271.1
$ws.Range("L46").Value = 487.25
$ws.Range("M46").Value = -1083.1
$ws.Range("N46").Value = -863.25
$ws.Range("H68").Value = 2646.3076
$ws.Range("I68").Value = 2889.111
$ws.Range("J68").Value = 2100
$ws.Range("K68").Value = 2889.111
$ws.Range("L68").Value = 2100
$ws.Range("M68").Value = -2140.111
$ws.Range("N68").Value = -3598
$ws.Range("H71").Value = 2646.3076
$ws.Range("I71").Value = 2889.111
$ws.Range("J71").Value = 2100
$ws.Range("K71").Value = 14445.555
$ws.Range("L71").Value = 10500
$ws.Range("M71").Value = -10701.555
$ws.Range("N71").Value = -17988
$ws.Range("H122").Value = 7901.1665
$ws.Range("I122").Value = 10107.708
$ws.Range("J122").Value = 3488.0833
$ws.Range("K122").Value = 30323.124
$ws.Range("L122").Value = 10464.2499
$ws.Range("M122").Value = -27873.124
$ws.Range("N122").Value = -15364.2499
$ws.Range("H132").Value = 1911.5641
$ws.Range("I132").Value = 1477.92
$ws.Range("J132").Value = 2685.9285
$ws.Range("K132").Value = 4433.76
$ws.Range("L132").Value = 8057.7855
$ws.Range("M132").Value = -1903.76
$ws.Range("N132").Value = -13117.7855
